# Regenerate the localization-status report for archive:
#  - flip the "Ready for handoff" status value to "In Translation"
#    everywhere it's used (Overview sheet + each per-locale sheet)
#  - re-autofit the "Status" column (and its mirror on the Overview
#    sheet) to the new, shorter text

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text on every sheet that shows it -----------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2. Re-autofit the affected "Status" columns to the new width -----
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
